$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 143
$ws.Cells.Item(143, 1).Value = 142
$ws.Cells.Item(143, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(143, 3).Value = "8:05 PM"
$ws.Cells.Item(143, 4).Value = "U23817"
$ws.Cells.Item(143, 5).Value = "Paris"
$ws.Cells.Item(143, 6).Value = "(CDG)"
$ws.Cells.Item(143, 7).Value = "easyJet "
$ws.Cells.Item(143, 8).Value = "A320"
$ws.Cells.Item(143, 9).Value = "(OE-IVD)"
$ws.Cells.Item(143, 10).Value = "8:04 PM"
$ws.Cells.Item(143, 12).Value = "0 hours, -1 minutes"

# Row 144
$ws.Cells.Item(144, 1).Value = 143
$ws.Cells.Item(144, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(144, 3).Value = "8:50 PM"
$ws.Cells.Item(144, 4).Value = "LO3921"
$ws.Cells.Item(144, 5).Value = "Warsaw"
$ws.Cells.Item(144, 6).Value = "(WAW)"
$ws.Cells.Item(144, 7).Value = "LOT "
$ws.Cells.Item(144, 8).Value = "E170"
$ws.Cells.Item(144, 9).Value = "(SP-LDF)"
$ws.Cells.Item(144, 10).Value = "8:42 PM"
$ws.Cells.Item(144, 12).Value = "0 hours, -8 minutes"

# Row 145
$ws.Cells.Item(145, 1).Value = 144
$ws.Cells.Item(145, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(145, 3).Value = "9:10 PM"
$ws.Cells.Item(145, 4).Value = "DY1042"
$ws.Cells.Item(145, 5).Value = "Oslo"
$ws.Cells.Item(145, 6).Value = "(OSL)"
$ws.Cells.Item(145, 7).Value = "Norwegian "
$ws.Cells.Item(145, 8).Value = "B738"
$ws.Cells.Item(145, 9).Value = "(LN-NIH)"
$ws.Cells.Item(145, 10).Value = "9:06 PM"
$ws.Cells.Item(145, 12).Value = "0 hours, -4 minutes"

# Row 146
$ws.Cells.Item(146, 1).Value = 145
$ws.Cells.Item(146, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(146, 3).Value = "9:15 PM"
$ws.Cells.Item(146, 4).Value = "FR2436"
$ws.Cells.Item(146, 5).Value = "London"
$ws.Cells.Item(146, 6).Value = "(STN)"
$ws.Cells.Item(146, 7).Value = "Ryanair "
$ws.Cells.Item(146, 8).Value = "B738"
$ws.Cells.Item(146, 9).Value = "(EI-EBY)"
$ws.Cells.Item(146, 10).Value = "8:51 PM"
$ws.Cells.Item(146, 12).Value = "0 hours, -24 minutes"

# Row 147
$ws.Cells.Item(147, 1).Value = 146
$ws.Cells.Item(147, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(147, 3).Value = "9:20 PM"
$ws.Cells.Item(147, 4).Value = "FR6265"
$ws.Cells.Item(147, 5).Value = "Poznan"
$ws.Cells.Item(147, 6).Value = "(POZ)"
$ws.Cells.Item(147, 7).Value = "Ryanair "
$ws.Cells.Item(147, 8).Value = "B738"
$ws.Cells.Item(147, 9).Value = "(SP-RSM)"
$ws.Cells.Item(147, 10).Value = "9:00 PM"
$ws.Cells.Item(147, 12).Value = "0 hours, -20 minutes"

# Row 148
$ws.Cells.Item(148, 1).Value = 147
$ws.Cells.Item(148, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(148, 3).Value = "9:40 PM"
$ws.Cells.Item(148, 4).Value = "FR6233"
$ws.Cells.Item(148, 5).Value = "Palermo"
$ws.Cells.Item(148, 6).Value = "(PMO)"
$ws.Cells.Item(148, 7).Value = "Ryanair "
$ws.Cells.Item(148, 8).Value = "B38M"
$ws.Cells.Item(148, 9).Value = "(SP-RZL)"
$ws.Cells.Item(148, 10).Value = "9:37 PM"
$ws.Cells.Item(148, 12).Value = "0 hours, -3 minutes"

# Row 149
$ws.Cells.Item(149, 1).Value = 148
$ws.Cells.Item(149, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(149, 3).Value = "9:55 PM"
$ws.Cells.Item(149, 4).Value = "U21871"
$ws.Cells.Item(149, 5).Value = "Manchester"
$ws.Cells.Item(149, 6).Value = "(MAN)"
$ws.Cells.Item(149, 7).Value = "easyJet "
$ws.Cells.Item(149, 8).Value = "A320"
$ws.Cells.Item(149, 9).Value = "(G-EZUS)"
$ws.Cells.Item(149, 10).Value = "9:44 PM"
$ws.Cells.Item(149, 12).Value = "0 hours, -11 minutes"

# Row 150
$ws.Cells.Item(150, 1).Value = 149
$ws.Cells.Item(150, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(150, 3).Value = "10:05 PM"
$ws.Cells.Item(150, 4).Value = "LH1624"
$ws.Cells.Item(150, 5).Value = "Munich"
$ws.Cells.Item(150, 6).Value = "(MUC)"
$ws.Cells.Item(150, 7).Value = "Lufthansa "
$ws.Cells.Item(150, 8).Value = "A320"
$ws.Cells.Item(150, 9).Value = "(D-AIZC)"
$ws.Cells.Item(150, 10).Value = "9:53 PM"
$ws.Cells.Item(150, 12).Value = "0 hours, -12 minutes"

# Row 151
$ws.Cells.Item(151, 1).Value = 150
$ws.Cells.Item(151, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(151, 3).Value = "10:10 PM"
$ws.Cells.Item(151, 4).Value = "FR4204"
$ws.Cells.Item(151, 5).Value = "Glasgow"
$ws.Cells.Item(151, 6).Value = "(GLA)"
$ws.Cells.Item(151, 7).Value = "Buzz "
$ws.Cells.Item(151, 8).Value = "B38M"
$ws.Cells.Item(151, 9).Value = "(SP-RZB)"
$ws.Cells.Item(151, 10).Value = "10:07 PM"
$ws.Cells.Item(151, 12).Value = "0 hours, -3 minutes"

# Row 152
$ws.Cells.Item(152, 1).Value = 151
$ws.Cells.Item(152, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(152, 3).Value = "10:20 PM"
$ws.Cells.Item(152, 4).Value = "FR1813"
$ws.Cells.Item(152, 5).Value = "London"
$ws.Cells.Item(152, 6).Value = "(LTN)"
$ws.Cells.Item(152, 7).Value = "Ryanair "
$ws.Cells.Item(152, 8).Value = "B738"
$ws.Cells.Item(152, 9).Value = "(SP-RKU)"
$ws.Cells.Item(152, 10).Value = "10:16 PM"
$ws.Cells.Item(152, 12).Value = "0 hours, -4 minutes"

# Row 153
$ws.Cells.Item(153, 1).Value = 152
$ws.Cells.Item(153, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(153, 3).Value = "10:40 PM"
$ws.Cells.Item(153, 4).Value = "FR3721"
$ws.Cells.Item(153, 5).Value = "Billund"
$ws.Cells.Item(153, 6).Value = "(BLL)"
$ws.Cells.Item(153, 7).Value = "Buzz "
$ws.Cells.Item(153, 8).Value = "B38M"
$ws.Cells.Item(153, 9).Value = "(SP-RZH)"
$ws.Cells.Item(153, 10).Value = "10:36 PM"
$ws.Cells.Item(153, 12).Value = "0 hours, -4 minutes"
